# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) previously held the "Strike#" count; this
# re-derives/rewrites that column with the real strikeout totals (K) per
# start, for rows 2-18 (the 17 most recent starts).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 5
    3  = 4
    4  = 2
    5  = 0
    6  = 3
    7  = 3
    8  = 6
    9  = 3
    10 = 4
    11 = 9
    12 = 9
    13 = 3
    14 = 8
    15 = 2
    16 = 2
    17 = 7
    18 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
